$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 819.2
$ws.Range("J32").Value = 819.2
$ws.Range("L32").Value = 819.2
$ws.Range("N32").Value = -1471.2
$ws.Range("H43").Value = 4178.273
$ws.Range("I43").Value = 4000.3333
$ws.Range("J43").Value = 4245
$ws.Range("K43").Value = 4000.3333
$ws.Range("L43").Value = 4245
$ws.Range("M43").Value = -3931.3333
$ws.Range("N43").Value = -4383
$ws.Range("H62").Value = 3519.7856
$ws.Range("I62").Value = 3752.5454
$ws.Range("K62").Value = 3752.5454
$ws.Range("M62").Value = -3128.5454
$ws.Range("H64").Value = 71432130
$ws.Range("I64").Value = 4045.875
$ws.Range("J64").Value = 166669570
$ws.Range("K64").Value = 4045.875
$ws.Range("L64").Value = 166669570
$ws.Range("M64").Value = -3797.875
$ws.Range("N64").Value = -166670066
$ws.Range("H65").Value = 3519.7856
$ws.Range("I65").Value = 3752.5454
$ws.Range("K65").Value = 18762.727
$ws.Range("M65").Value = -15642.727
$ws.Range("H67").Value = 71432130
$ws.Range("I67").Value = 4045.875
$ws.Range("J67").Value = 166669570
$ws.Range("K67").Value = 4045.875
$ws.Range("L67").Value = 166669570
$ws.Range("M67").Value = -3187.875
$ws.Range("N67").Value = -166671286
$ws.Range("H74").Value = 10666.518
$ws.Range("I74").Value = 11174.192
$ws.Range("J74").Value = 6266.6665
$ws.Range("K74").Value = 11174.192
$ws.Range("L74").Value = 6266.6665
$ws.Range("M74").Value = -10238.192
$ws.Range("N74").Value = -8138.6665
$ws.Range("H77").Value = 10666.518
$ws.Range("I77").Value = 11174.192
$ws.Range("J77").Value = 6266.6665
$ws.Range("K77").Value = 55870.95999999999
$ws.Range("L77").Value = 31333.3325
$ws.Range("M77").Value = -51190.95999999999
$ws.Range("N77").Value = -40693.3325
$ws.Range("H132").Value = 4880.592
$ws.Range("I132").Value = 5283.8096
$ws.Range("K132").Value = 15851.4288
$ws.Range("M132").Value = -13321.4288
$ws.Range("H137").Value = 1192944.9
$ws.Range("I137").Value = 1668691.6
$ws.Range("J137").Value = 3578.25
$ws.Range("K137").Value = 5006074.800000001
$ws.Range("L137").Value = 10734.75
$ws.Range("M137").Value = -5003524.800000001
$ws.Range("N137").Value = -15834.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2983.8545
$ws.Range("I32").Value = 2983.8545
$ws.Range("K32").Value = 2983.8545
$ws.Range("M32").Value = -2696.8545
$ws.Range("H45").Value = 15495.875
$ws.Range("I45").Value = 18915
$ws.Range("J45").Value = 3284.7144
$ws.Range("K45").Value = 18915
$ws.Range("L45").Value = 3284.7144
$ws.Range("M45").Value = -18538
$ws.Range("N45").Value = -4038.7144
$ws.Range("H63").Value = 950.5
$ws.Range("I63").Value = 920.6
$ws.Range("J63").Value = 1100
$ws.Range("K63").Value = 920.6
$ws.Range("L63").Value = 1100
$ws.Range("M63").Value = -234.6
$ws.Range("N63").Value = -2472
$ws.Range("H66").Value = 950.5
$ws.Range("I66").Value = 920.6
$ws.Range("J66").Value = 1100
$ws.Range("K66").Value = 4603
$ws.Range("L66").Value = 5500
$ws.Range("M66").Value = -1171
$ws.Range("N66").Value = -12364
$ws.Range("H74").Value = 293912.3
$ws.Range("I74").Value = 348082.8
$ws.Range("K74").Value = 348082.8
$ws.Range("M74").Value = -347208.8
$ws.Range("H77").Value = 293912.3
$ws.Range("I77").Value = 348082.8
$ws.Range("K77").Value = 1740414
$ws.Range("M77").Value = -1736046
$ws.Range("H102").Value = 5670
$ws.Range("I102").Value = 5038.4
$ws.Range("K102").Value = 5038.4
$ws.Range("M102").Value = -3416.4
$ws.Range("H132").Value = 6291942
$ws.Range("I132").Value = 2262.2632
$ws.Range("K132").Value = 6786.7896
$ws.Range("M132").Value = -4256.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 142859000
$ws.Range("I94").Value = 285714880
$ws.Range("J94").Value = 3156.2856
$ws.Range("K94").Value = 285714880
$ws.Range("L94").Value = 3156.2856
$ws.Range("M94").Value = -285714429
$ws.Range("N94").Value = -4058.2856
$ws.Range("H107").Value = 1189.409
$ws.Range("I107").Value = 1134
$ws.Range("J107").Value = 1540.3334
$ws.Range("K107").Value = 1134
$ws.Range("L107").Value = 1540.3334
$ws.Range("M107").Value = 786
$ws.Range("N107").Value = -5380.3334
$ws.Range("H134").Value = 2614.745
$ws.Range("J134").Value = 4375
$ws.Range("L134").Value = 13125
$ws.Range("N134").Value = -18195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5023.339
$ws.Range("I31").Value = 4090.848
$ws.Range("K31").Value = 4090.848
$ws.Range("M31").Value = -3795.848
$ws.Range("H34").Value = 5023.339
$ws.Range("I34").Value = 4090.848
$ws.Range("K34").Value = 4090.848
$ws.Range("M34").Value = -3888.848
$ws.Range("H99").Value = 3600
$ws.Range("I99").Value = 3416.6667
$ws.Range("K99").Value = 3416.6667
$ws.Range("M99").Value = -1918.6667
$ws.Range("H122").Value = 2941.1904
$ws.Range("I122").Value = 2696.6
$ws.Range("J122").Value = 3552.6667
$ws.Range("K122").Value = 8089.799999999999
$ws.Range("L122").Value = 10658.0001
$ws.Range("M122").Value = -5639.799999999999
$ws.Range("N122").Value = -15558.0001
$ws.Range("H126").Value = 3600
$ws.Range("I126").Value = 3416.6667
$ws.Range("K126").Value = 10250.0001
$ws.Range("M126").Value = -7780.000100000001
$ws.Range("H134").Value = 2743.2307
$ws.Range("I134").Value = 2666.0881
$ws.Range("K134").Value = 7998.2643
$ws.Range("M134").Value = -5463.2643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 400
$ws.Range("I70").Value = 400
$ws.Range("K70").Value = 1200
$ws.Range("M70").Value = -885
$ws.Range("H73").Value = 400
$ws.Range("I73").Value = 400
$ws.Range("K73").Value = 1200
$ws.Range("M73").Value = -108
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H140").Value = 10718.828
$ws.Range("I140").Value = 7186.7827
$ws.Range("J140").Value = 17488.584
$ws.Range("K140").Value = 21560.3481
$ws.Range("L140").Value = 52465.75199999999
$ws.Range("M140").Value = -16380.3481
$ws.Range("N140").Value = -62825.75199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5203.963
$ws.Range("I102").Value = 1595.6086
$ws.Range("J102").Value = 25952
$ws.Range("K102").Value = 1595.6086
$ws.Range("L102").Value = 25952
$ws.Range("M102").Value = 26.39139999999998
$ws.Range("N102").Value = -29196
$ws.Range("H122").Value = 3707.8462
$ws.Range("I122").Value = 3084.35
$ws.Range("J122").Value = 5786.1665
$ws.Range("K122").Value = 9253.049999999999
$ws.Range("L122").Value = 17358.4995
$ws.Range("M122").Value = -6803.049999999999
$ws.Range("N122").Value = -22258.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4561.385
$ws.Range("I7").Value = 3771.818
$ws.Range("J7").Value = 8904
$ws.Range("K7").Value = 3771.818
$ws.Range("L7").Value = 8904
$ws.Range("M7").Value = -3659.818
$ws.Range("N7").Value = -9128
$ws.Range("H61").Value = 29409.75
$ws.Range("I61").Value = 2554.2
$ws.Range("J61").Value = 74169
$ws.Range("K61").Value = 2554.2
$ws.Range("L61").Value = 74169
$ws.Range("M61").Value = -2352.2
$ws.Range("N61").Value = -74573
$ws.Range("H93").Value = 460.75
$ws.Range("I93").Value = 475.7143
$ws.Range("J93").Value = 439.8
$ws.Range("K93").Value = 475.7143
$ws.Range("L93").Value = 439.8
$ws.Range("M93").Value = 772.2857
$ws.Range("N93").Value = -2935.8
$ws.Range("H113").Value = 29409.75
$ws.Range("I113").Value = 2554.2
$ws.Range("J113").Value = 74169
$ws.Range("K113").Value = 2554.2
$ws.Range("L113").Value = 74169
$ws.Range("M113").Value = -384.1999999999998
$ws.Range("N113").Value = -78509
$ws.Range("H126").Value = 4561.385
$ws.Range("I126").Value = 3771.818
$ws.Range("J126").Value = 8904
$ws.Range("K126").Value = 11315.454
$ws.Range("L126").Value = 26712
$ws.Range("M126").Value = -8845.454000000002
$ws.Range("N126").Value = -31652
$ws.Range("H132").Value = 2971.0513
$ws.Range("I132").Value = 2376.3333
$ws.Range("K132").Value = 7128.999899999999
$ws.Range("M132").Value = -4598.999899999999
$ws.Range("H136").Value = 2760.7273
$ws.Range("I136").Value = 2681.6667
$ws.Range("J136").Value = 3116.5
$ws.Range("K136").Value = 8045.000100000001
$ws.Range("L136").Value = 9349.5
$ws.Range("M136").Value = -5495.000100000001
$ws.Range("N136").Value = -14449.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7754953.5
$ws.Range("I132").Value = 10419432
$ws.Range("J132").Value = 3744.818
$ws.Range("K132").Value = 31258296
$ws.Range("L132").Value = 11234.454
$ws.Range("M132").Value = -31255766
$ws.Range("N132").Value = -16294.454
$ws.Range("H136").Value = 6667.2163
$ws.Range("I136").Value = 5861.032
$ws.Range("J136").Value = 10832.5
$ws.Range("K136").Value = 17583.096
$ws.Range("L136").Value = 32497.5
$ws.Range("M136").Value = -15033.096
$ws.Range("N136").Value = -37597.5
